$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Rows 17 and 18: ShibaInu and BinanceUSD swapped positions with updated price/volume/link
Set-TextValue "B17" "ShibaInu"
Set-TextValue "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.00001141"
Set-TextValue "E17" "  +1.13%  "

Set-TextValue "B18" "BinanceUSD"
Set-TextValue "C18" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D18" "1.005"
Set-TextValue "E18" "  +0.23%  "

# Remaining rows: update Price (D) and/or Volume(1h) (E) values as scraped
Set-TextValue "D2" "30.421.55"
Set-TextValue "E2" "  -0.24%  "
Set-TextValue "D3" "2.100.23"
Set-TextValue "E3" "  -0.33%  "
Set-TextValue "D4" "1.004"
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "334.22"
Set-TextValue "E5" "  +1.50%  "
Set-TextValue "D6" "1.004"
Set-TextValue "E6" "  +0.30%  "
Set-TextValue "E7" "  -1.03%  "
Set-TextValue "D8" "0.4526"
Set-TextValue "D9" "53.93"
Set-TextValue "E9" "  +14.08%  "
Set-TextValue "D10" "0.08865"
Set-TextValue "E10" "  -0.09%  "
Set-TextValue "E11" "  +1.36%  "
Set-TextValue "D12" "24.04"
Set-TextValue "E12" "  -1.98%  "
Set-TextValue "D13" "2.093.57"
Set-TextValue "E13" "  -0.77%  "
Set-TextValue "D14" "6.790"
Set-TextValue "E14" "  +0.94%  "
Set-TextValue "D15" "8.006"
Set-TextValue "E15" "  +3.15%  "
Set-TextValue "D16" "96.68"
Set-TextValue "E16" "  +0.35%  "
Set-TextValue "D19" "0.06629"
Set-TextValue "E19" "  -0.22%  "
Set-TextValue "D20" "19.15"
Set-TextValue "E20" "  +0.64%  "
Set-TextValue "D21" "1.004"
Set-TextValue "E21" "  +0.28%  "
Set-TextValue "D22" "6.264"
Set-TextValue "E22" "  -0.94%  "
Set-TextValue "D23" "30.486.81"
Set-TextValue "E23" "  -0.21%  "
Set-TextValue "E24" "  +0.07%  "
Set-TextValue "D25" "2.337"
Set-TextValue "E25" "  -0.79%  "
Set-TextValue "D26" "2.341.03"
Set-TextValue "E26" "  -0.56%  "
Set-TextValue "D27" "22.12"
Set-TextValue "E27" "  -1.38%  "
Set-TextValue "D28" "162.51"
Set-TextValue "E28" "  +0.26%  "
Set-TextValue "E29" "  -2.84%  "
Set-TextValue "D30" "132.85"
Set-TextValue "E30" "  +0.13%  "
Set-TextValue "D31" "1.199"
Set-TextValue "E31" "  -0.92%  "
Set-TextValue "E32" "  -1.11%  "
Set-TextValue "D33" "1.648"
Set-TextValue "E33" "  -1.18%  "
Set-TextValue "E34" "  +2.62%  "
Set-TextValue "D35" "3.949"
Set-TextValue "E35" "  +0.65%  "
Set-TextValue "D36" "10.36"
Set-TextValue "E36" "  +3.61%  "
Set-TextValue "D37" "5.815"
Set-TextValue "E37" "  +5.96%  "
Set-TextValue "D38" "0.02570"
Set-TextValue "E38" "  -0.63%  "
Set-TextValue "D39" "0.06833"
Set-TextValue "E39" "  +2.24%  "
Set-TextValue "E40" "  +0.28%  "
Set-TextValue "D41" "12.69"
Set-TextValue "E41" "  -0.09%  "
Set-TextValue "D42" "0.6853"
Set-TextValue "E42" "  +0.84%  "
Set-TextValue "D43" "1.241"
Set-TextValue "E43" "  -2.49%  "
Set-TextValue "D44" "2.319"
Set-TextValue "E44" "  +4.98%  "
Set-TextValue "D45" "13.98"
Set-TextValue "E45" "  -0.38%  "
Set-TextValue "D46" "0.6333"
Set-TextValue "E46" "  -0.77%  "
Set-TextValue "E47" "  +1.12%  "
Set-TextValue "D48" "1.245"
Set-TextValue "E48" "  -0.48%  "
Set-TextValue "D49" "0.00000000346"
Set-TextValue "E49" "  +17.68%  "
Set-TextValue "E50" "  +0.30%  "
Set-TextValue "E51" "  +0.16%  "
